$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false | Out-Null

# 1. Remove the "capacityFactor" sheet entirely; its data is replaced by the
#    new "flow" sheet (added below) with a reshaped layout.
$capacityFactor = $wb.Worksheets.Item("capacityFactor")
$capacityFactor.Delete() | Out-Null

# 2. Rename "p_unit" to "unit_p" (content is unchanged).
$pUnit = $wb.Worksheets.Item("p_unit")
$pUnit.Name = "unit_p"

# 3. Reshape the "unit_ts" sheet data.
$unitTs = $wb.Worksheets.Item("unit_ts")
$unitTs.Cells.Clear()
$unitTs.Range("D1").Value = "gas_turbine"
$unitTs.Range("A2").Value = "Base"
$unitTs.Range("B2").Value = "eff01_ts"
$unitTs.Range("C2").Value = "t000001"
$unitTs.Range("D2").Value = 0.5

# 4. Add a new "flow" sheet at the end of the workbook with the reshaped
#    capacity-factor data. Copy an existing sheet first so the new sheet
#    keeps the same sheet-level formatting conventions used throughout
#    this workbook, then clear it out and fill in the new values.
$templateSheet = $wb.Worksheets.Item("reserveDemand")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$templateSheet.Copy($null, $lastSheet)
$flow = $wb.Worksheets.Item($wb.Worksheets.Count)
$flow.Name = "flow"
$flow.Cells.Clear()

$flow.Range("D1").Value = "wind1"

$flow.Range("A2").Value = "capacityFactor"
$flow.Range("B2").Value = "Base"
$flow.Range("C2").Value = "t000001"
$flow.Range("D2").Value = 0.4

$flow.Range("A3").Value = "capacityFactor"
$flow.Range("B3").Value = "Base"
$flow.Range("C3").Value = "t000002"
$flow.Range("D3").Value = 0.4

$flow.Range("A4").Value = "capacityFactor"
$flow.Range("B4").Value = "Base"
$flow.Range("C4").Value = "t000003"
$flow.Range("D4").Value = 0.4

$flow.Range("A5").Value = "capacityFactor"
$flow.Range("B5").Value = "Base"
$flow.Range("C5").Value = "t000004"
$flow.Range("D5").Value = 0.4

$flow.Range("A6").Value = "capacityFactor"
$flow.Range("B6").Value = "Base"
$flow.Range("C6").Value = "t000005"
$flow.Range("D6").Value = 0.4
